$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.436.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.30%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.189.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.02%  '

# Row 4
$ws.Range("E4").Value = '  -0.21%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.53%  '

# Row 6
$ws.Range("E6").Value = '  -1.67%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '67.27'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.73%  '

# Row 8
$ws.Range("E8").Value = '  -0.13%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.14%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.89'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.67%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.91%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0934'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.74%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.44%  '

# Row 14
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.103'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.73%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.515.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.10%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.188.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.355.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0951'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.28%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.66%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.69%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.04%  '

# Row 25
$ws.Range("E25").Value = '  -3.49%  '

# Row 26
$ws.Range("E26").Value = '  +0.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.21%  '

# Row 28
$ws.Range("E28").Value = '  -5.38%  '

# Row 29
$ws.Range("E29").Value = '  -2.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.03%  '

# Row 31
$ws.Range("E31").Value = '  -6.58%  '

# Row 32
$ws.Range("E32").Value = '  -2.24%  '

# Row 33
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.120'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.04%  '

# Row 34
$ws.Range("E34").Value = '  +3.06%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0783'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.122'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.76%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.57%  '

# Row 38
$ws.Range("E38").Value = '  -2.53%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.81%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0305'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.35%  '

# Row 41
$ws.Range("E41").Value = '  -2.49%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.14'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.72%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.24%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.41%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.06%  '

# Row 46
$ws.Range("E46").Value = '  -3.82%  '

# Row 47
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.49%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1000'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.82%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.52%  '

# Row 50
$ws.Range("E50").Value = '  -1.68%  '

# Row 51
$ws.Range("E51").Value = '  +5.73%  '
